$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the now-unused rows 4 and 5 first (A4/B4 = "Appel", A5/B5 = "Cookie")
$ws.Rows("4:5").Delete()

# Insert new "username" header in B1, shifting existing B1 content to C1
$ws.Range("B1").Value2 = "username"

# C1 gets the "massage one account" header, copying B1's style (bold, bordered, centered)
$ws.Range("C1").Value2 = "massage one account"
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)

# Row 2: update B2, add C2
$ws.Range("B2").Value2 = "SabaMosaybie"
$ws.Range("C2").Value2 = "hello"

# Row 3: update B3, add C3
$ws.Range("B3").Value2 = "narges__pv310"
$ws.Range("C3").Value2 = "my name is selenium py"
